# 6th nov test cases
# Adds three new test-case rows (FW_UI_0002, FW_UI_0003, FW_UI_0004) to the
# "Test Cases" table, matching formatting from the existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by three rows (it currently has one data row -> four).
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Copy the formatting of the existing data row down onto the new rows so the
# borders/fonts/alignment used for the table body stay consistent.
$formatSrc = $ws.Range("A2:F2")
$formatSrc.Copy($ws.Range("A3:F3"))
$formatSrc.Copy($ws.Range("A4:F4"))
$formatSrc.Copy($ws.Range("A5:F5"))
$excel.CutCopyMode = $false

# Match the row height used throughout the rest of the table.
$ws.Rows.Item(3).RowHeight = 13.5
$ws.Rows.Item(4).RowHeight = 13.5
$ws.Rows.Item(5).RowHeight = 13.5

# Column order: RunTest, TC_ID, ScriptName, Parameters, Class Name, Description
$rows = @(
  @("No",  "FW_UI_0001", "VerifyPublishScenarios",  "Story,TestStoryHeadline,Test,Test,TestMe,Test.V,Test,T,Taslic/,TestStory", "PublishCases",  "Aim of the script is to Create a Multiple New Stories  and Publish the same"),
  @("No",  "FW_UI_0002", "VerifyPublishScenarios",  "Alert,TestAlertHeadline,Test,Test,TestMe,Test.V, Test,T,Taslic/,NA",       "PublishCases",  "Aim of the script is to Create a  Multiple New Alerts,Publish and Verify in the GRID"),
  @("No",  "FW_UI_0003", "VerifyPublishScenarios",  "Econ,64424509456,7,5,NA,NA, NA,NA,NA,NA",                                  "PublishCases",  "Aim of the script is to Create Mulitple New  Econs , Publilsh and Verify in the Grid"),
  @("Yes", "FW_UI_0004", "VerifyTemplateScenarios", "Story,TestStoryHeadline,Test,Test,TestMe,Test.V,Test,T,Taslic/,TestStory", "TemplateCases", "Aim of the script is to Create a New Story and Save as Template")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = 2 + $i
  $rowData = $rows[$i]
  for ($c = 0; $c -lt $rowData.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value2 = $rowData[$c]
  }
}

$ws.Range("D5").Select() | Out-Null
